# Update stimulus presentation time-logging: refresh timestamp-based file
# names / sheet names and fix the RS_TO eyes-open/eyes-closed ordering.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO -----------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Name = "GNG_TO-16512555350024083"
$ws.Range("B2").Value = "go_stims-16512555349584014.csv"
$ws.Range("B3").Value = "GNG_stims-16512555349853992.csv"
$ws.Range("B4").Value = "go_stims-16512555349874117.csv"
$ws.Range("B5").Value = "GNG_stims-16512555350014086.csv"

# --- Sheet 2: NB_TO --------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Name = "NB_TO-16512555370722148"
$ws.Range("B2").Value = "ZB-match_4-1651255535336214.csv"
$ws.Range("B3").Value = "TB-16512555362112145.csv"
$ws.Range("B4").Value = "ZB-match_1-16512555350944104.csv"
$ws.Range("B5").Value = "OB-16512555358282135.csv"
$ws.Range("B6").Value = "TB-16512555370502155.csv"
$ws.Range("B7").Value = "OB-1651255536184218.csv"
$ws.Range("B8").Value = "ZB-match_8-1651255535029402.csv"
$ws.Range("B9").Value = "TB-16512555363792148.csv"
$ws.Range("B10").Value = "OB-1651255535609213.csv"

# --- Sheet 3: RS_TO ---------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Name = "RS_TO-1651255537078218"
$ws.Range("B2").Value = "eyes open"
$ws.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Name = "TOL_TO-1651255537137214"
$ws.Range("B2").Value = "MM_stims-16512555371042135.csv"
$ws.Range("B3").Value = "ZM_stims-16512555370802197.csv"
$ws.Range("B4").Value = "MM_stims-1651255537120215.csv"
$ws.Range("B5").Value = "ZM_stims-165125553710522.csv"
$ws.Range("B6").Value = "MM_stims-16512555371362162.csv"
$ws.Range("B7").Value = "ZM_stims-16512555371212165.csv"

# --- Sheet 5: vSAT_TO ---------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Name = "vSAT_TO-16512555372162137"
$ws.Range("B2").Value = "vSAT_stims-1651255537184214.csv"
$ws.Range("B3").Value = "SAT_stims-165125553714122.csv"
$ws.Range("B4").Value = "SAT_stims-1651255537168215.csv"
$ws.Range("B5").Value = "vSAT_stims-1651255537200213.csv"
